# The commit regenerates the "K" column (strikeouts, column G) in the
# save_data sheet for hoffman_jeff with freshly (re)computed values.
# Apply the new computed values to column G, rows 2-67.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 1
    4  = 3
    5  = 1
    6  = 2
    7  = 0
    8  = 1
    10 = 1
    11 = 0
    12 = 1
    15 = 1
    16 = 2
    17 = 0
    18 = 2
    19 = 2
    20 = 2
    21 = 1
    22 = 2
    23 = 2
    24 = 1
    25 = 0
    26 = 2
    27 = 1
    28 = 0
    29 = 2
    30 = 0
    31 = 2
    32 = 2
    33 = 1
    34 = 1
    35 = 2
    36 = 0
    37 = 0
    38 = 1
    39 = 2
    40 = 2
    41 = 0
    42 = 1
    43 = 1
    44 = 2
    45 = 2
    46 = 2
    47 = 2
    48 = 2
    49 = 2
    50 = 1
    51 = 1
    52 = 2
    53 = 0
    54 = 2
    55 = 0
    56 = 1
    57 = 0
    58 = 1
    59 = 1
    60 = 3
    61 = 4
    62 = 2
    63 = 2
    64 = 2
    65 = 2
    66 = 2
    67 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
